$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 0. Remove the existing "_GoBack" bookmark from the last paragraph;
#    it will be re-created at the end of paragraph 5 further down.
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ------------------------------------------------------------------
# Helper positions: paragraph 5 is the one ending in
# "http://i.giphy.com/xT0BKzyF3d7ljlJpug.gif" (inside a hyperlink).
# We always re-fetch the paragraph + insertion point fresh after each
# edit because earlier edits shift later character offsets.
# ------------------------------------------------------------------

function Get-P5End {
    $p = $d.Paragraphs(5)
    return $p.Range.End - 1
}

# ------------------------------------------------------------------
# 1. Two line breaks, styled like the Hyperlink run style.
# ------------------------------------------------------------------
$pos = Get-P5End
$ip = $d.Range($pos, $pos)
$ip.InsertBreak(6)
$pos2 = Get-P5End
$rStyle = $d.Range($pos2 - 1, $pos2)
$rStyle.Style = "Hyperlink"

$pos = Get-P5End
$ip = $d.Range($pos, $pos)
$ip.InsertBreak(6)
$pos2 = Get-P5End
$rStyle = $d.Range($pos2 - 1, $pos2)
$rStyle.Style = "Hyperlink"

# ------------------------------------------------------------------
# 2. "Patterny hexy thingy:" text in red, Hyperlink style, no underline.
# ------------------------------------------------------------------
$text = "Patterny hexy thingy:"
$pos = Get-P5End
$ip = $d.Range($pos, $pos)
$ip.InsertAfter($text)
$pos2 = Get-P5End
$rText = $d.Range($pos2 - $text.Length, $pos2)
$rText.Style = "Hyperlink"
$rText.Font.Color = 255
$rText.Font.Underline = 0

# ------------------------------------------------------------------
# 3. Line break, styled like Hyperlink.
# ------------------------------------------------------------------
$pos = Get-P5End
$ip = $d.Range($pos, $pos)
$ip.InsertBreak(6)
$pos2 = Get-P5End
$rStyle = $d.Range($pos2 - 1, $pos2)
$rStyle.Style = "Hyperlink"

# ------------------------------------------------------------------
# 4. Hyperlink to giphy media URL.
# ------------------------------------------------------------------
$url1 = "https://media.giphy.com/media/5xtDarE06fr4xurDL5m/giphy.gif"
$pos = Get-P5End
$ip = $d.Range($pos, $pos)
$ip.InsertAfter($url1)
$pos2 = Get-P5End
$rLink = $d.Range($pos2 - $url1.Length, $pos2)
$d.Hyperlinks.Add($rLink, $url1) | Out-Null

# ------------------------------------------------------------------
# 5. Plain line break (no style).
# ------------------------------------------------------------------
$pos = Get-P5End
$ip = $d.Range($pos, $pos)
$ip.InsertBreak(6)

# ------------------------------------------------------------------
# 6. Line break + "Trippy Tron Stairs Thingy:" plain text (same run).
# ------------------------------------------------------------------
$pos = Get-P5End
$ip = $d.Range($pos, $pos)
$ip.InsertBreak(6)
$ip2pos = Get-P5End
$ip2 = $d.Range($ip2pos, $ip2pos)
$ip2.InsertAfter("Trippy Tron Stairs Thingy:")

# ------------------------------------------------------------------
# 7. Plain line break (no style).
# ------------------------------------------------------------------
$pos = Get-P5End
$ip = $d.Range($pos, $pos)
$ip.InsertBreak(6)

# ------------------------------------------------------------------
# 8. Hyperlink to tumblr URL.
# ------------------------------------------------------------------
$url2 = "https://45.media.tumblr.com/64f3ed6cf9ed778515d3102b42583b40/tumblr_o3ac7aExOH1rsdpaso1_500.gif"
$pos = Get-P5End
$ip = $d.Range($pos, $pos)
$ip.InsertAfter($url2)
$pos2 = Get-P5End
$rLink2 = $d.Range($pos2 - $url2.Length, $pos2)
$d.Hyperlinks.Add($rLink2, $url2) | Out-Null

# ------------------------------------------------------------------
# 9. Two plain line breaks (no style).
# ------------------------------------------------------------------
$pos = Get-P5End
$ip = $d.Range($pos, $pos)
$ip.InsertBreak(6)

$pos = Get-P5End
$ip = $d.Range($pos, $pos)
$ip.InsertBreak(6)

# ------------------------------------------------------------------
# 10. Re-create the "_GoBack" bookmark at the end of paragraph 5.
# ------------------------------------------------------------------
$pos = Get-P5End
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output "Edit complete"
